$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B25" = 1.605796365659359;  "D25" = 2.418585580357664
    "B26" = 1.595170373197115;  "D26" = 2.423189895590971
    "B31" = 1.559211995598527;  "D31" = 2.259060943174633
    "B32" = 1.567714000147012;  "D32" = 2.289265470940818
    "B33" = 1.564079939769674;  "D33" = 2.264544154068898
    "B34" = 1.554941631335987;  "D34" = 2.187896682526145
    "B35" = 1.550489742549896;  "D35" = 2.239272596246528
    "B38" = 1.427985993580319;  "D38" = 1.945719472026729
    "B39" = 1.410548045123722;  "D39" = 1.882167269668976
    "B40" = 1.405378496021678;  "D40" = 1.853692822352101
    "B41" = 1.408893757585069;  "D41" = 1.834384170021593
    "B42" = 1.380351075228878;  "D42" = 1.815920663940857
    "B43" = 1.363977034209024;  "D43" = 1.778751645324429
    "B44" = 1.380901000651561;  "D44" = 1.766531543861792
    "B47" = 1.342218401785662;  "D47" = 1.688706241602388
    "B48" = 1.331246744782096;  "D48" = 1.669934466960038
    "B49" = 1.315140810373198;  "D49" = 1.639036470100827
    "B50" = 1.304783119122762;  "D50" = 1.625178676835295
    "B51" = 1.302501009379456;  "D51" = 1.598601300647098
    "B52" = 1.296052864122837;  "D52" = 1.588232197226381
    "B53" = 1.276281382275522;  "D53" = 1.549199684641708
    "B54" = 1.190921954421548;  "D54" = 1.825434111347368
    "B55" = 1.150148931499757;  "D55" = 1.699512378614168
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
